# Rename the "index" column to "i" and re-base its values to start at 0
# (was 1-based: 1..502, now 0-based: 0..501). Also narrow column A to fit
# the shorter header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header: rename column from "index" to "i" (this also renames the
# corresponding ListObject/table column automatically).
$ws.Range("A1").Value = "i"

# Data rows 2..503 held 1..502; shift every value down by one so the
# index column is zero-based.
for ($r = 2; $r -le 503; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Column A no longer needs to fit "index" (6 chars wide) - narrow it to
# fit "i" (width 4 in the saved OOXML).
$ws.Columns.Item(1).ColumnWidth = 3.17
